$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Actualiza base de datos EC ---

# "VALOR MORA" total (E11)
$ws.Range("E11").Value = 50434

# "Cant. Trabajadores" (C13) and "Cant. Periodos" (F13)
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1

# Remove the worker row for WILLIAM ISAAC ESCAMILLA VALENCIA (row 17 —
# CC 8765295, periodo 2111). The rows below (Felix's data row and the
# signature block) shift up to fill the gap.
$ws.Rows("17").Delete()

# "Nombre Trabajador" column best-fit width narrows now that the longest
# remaining worker name is shorter than the removed one.
$ws.Columns("D").ColumnWidth = 33.08
